$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh of the price data: the rows were reshuffled/updated with
# the latest week's figures. Write the new values for rows 2-7 (columns
# D, I, J, K, L, M, N, P, Q) to match the updated source export.

$rows = @(
    @{ Row = 2;  D = 44315; I = "Especial"; J = 10; K = 30000; L = 30000; M = 30000; N = "`$/caja 20 kilos empedrada"; P = 1500; Q = 20 },
    @{ Row = 3;  D = 44315; I = "Primera";  J = 20; K = 15000; L = 15000; M = 15000; N = "`$/caja 15 kilos granel";    P = 1000; Q = 15 },
    @{ Row = 4;  D = 44285; I = "Primera";  J = 20; K = 25000; L = 25000; M = 25000; N = "`$/caja 18 kilos empedrada"; P = 1389; Q = 18 },
    @{ Row = 5;  D = 44313; I = "Primera";  J = 20; K = 15000; L = 15000; M = 15000; N = "`$/caja 15 kilos empedrada"; P = 1000; Q = 15 },
    @{ Row = 6;  D = 44313; I = "Primera";  J = 20; K = 30000; L = 30000; M = 30000; N = "`$/caja 20 kilos empedrada"; P = 1500; Q = 20 },
    @{ Row = 7;  D = 44280; I = "Primera";  J = 30; K = 25000; L = 25000; M = 25000; N = "`$/caja 18 kilos empedrada"; P = 1389; Q = 18 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value  = $r.D   # D: Fecha
    $ws.Cells.Item($row, 9).Value  = $r.I   # I: Calidad
    $ws.Cells.Item($row, 10).Value = $r.J   # J: Volumen
    $ws.Cells.Item($row, 11).Value = $r.K   # K: Precio minimo
    $ws.Cells.Item($row, 12).Value = $r.L   # L: Precio maximo
    $ws.Cells.Item($row, 13).Value = $r.M   # M: Precio promedio ponderado
    $ws.Cells.Item($row, 14).Value = $r.N   # N: Unidad de comercializacion
    $ws.Cells.Item($row, 16).Value = $r.P   # P: Precio $/Kg
    $ws.Cells.Item($row, 17).Value = $r.Q   # Q: Kg o Unidades
}
